# Duplicate the current first sheet ("20190624") and place the copy
# before it, renaming the copy to "20190703". This mirrors the author's
# workflow of starting a new week's status sheet from a copy of the
# previous one, while the previous sheet ("20190624") keeps its data
# unchanged but is no longer the active/selected tab.

$wb = $excel.ActiveWorkbook

# The sheet to duplicate is the current first/active worksheet ("20190624").
$sourceSheet = $wb.Worksheets.Item(1)
$sourceSheetName = $sourceSheet.Name

# Copy it to a position right before itself -> the copy becomes the new
# first sheet, and the original slides down to position 2.
$sourceSheet.Copy($sourceSheet)

# The freshly-created copy is now the first sheet; rename it.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "20190703"

# The original sheet (still holding its original name/content) is now
# the second sheet.
$oldSheet = $wb.Worksheets.Item($sourceSheetName)

# Update the remembered selection on each tab: the new (active) sheet
# keeps the cursor at C48, while the original sheet's cursor moves to B39.
$newSheet.Activate()
$newSheet.Range("C48").Select()

$oldSheet.Activate()
$oldSheet.Range("B39").Select()

# Leave the new, duplicated sheet as the active tab.
$newSheet.Activate()
